$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "71.437.56"
$ws.Range("E2").Value = "  -0.25%  "

# Row 3
$ws.Range("D3").Value = "3.841.59"
$ws.Range("E3").Value = "  +0.68%  "

# Row 4
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "708.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.47%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.88%  "

# Row 7
$ws.Range("D7").Value = "3.839.21"
$ws.Range("E7").Value = "  +0.66%  "

# Row 8
$ws.Range("E8").Value = "  +0.07%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.527"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.48%  "

# Row 10
$ws.Range("E10").Value = "  -0.10%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.35"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.60%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.461"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.21%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000257"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.81%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.06"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.65%  "

# Row 15
$ws.Range("D15").Value = "4.627.95"
$ws.Range("E15").Value = "  +3.85%  "

# Row 16
$ws.Range("D16").Value = "3.796.70"
$ws.Range("E16").Value = "  -0.78%  "

# Row 17
$ws.Range("D17").Value = "71.473.36"
$ws.Range("E17").Value = "  -0.09%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.27"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.65%  "

# Row 19
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.51"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.37%  "

# Row 20
$ws.Range("B20").Value = "TRON"
$ws.Range("C20").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.115"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.37%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "499.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.24%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.67%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.735"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.67%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.55"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.13%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000146"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.72%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.59%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.67%  "

# Row 28
$ws.Range("D28").Value = "3.995.16"

# Row 29
$ws.Range("E29").Value = "  -2.33%  "

# Row 30
$ws.Range("E30").Value = "  -0.06%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.75%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.51"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.69%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.25"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.49%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.50"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.43%  "

# Row 35
$ws.Range("E35").Value = "  -4.96%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.26"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.40%  "

# Row 37
$ws.Range("D37").Value = "3.805.99"

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.998"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.31%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.104"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.09%  "

# Row 40
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.35"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.96%  "

# Row 41
$ws.Range("B41").Value = "Mantle"
$ws.Range("C41").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.04"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.05%  "

# Row 42
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.39"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.54%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.04"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.13%  "

# Row 44
$ws.Range("E44").Value = "  +0.00%  "

# Row 45
$ws.Range("E45").Value = "  +0.31%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.000318"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.47%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "164.16"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.54%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "431.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.53%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "49.05"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.44%  "

# Row 50
$ws.Range("E50").Value = "  +1.33%  "

# Row 51
$ws.Range("E51").Value = "  -0.35%  "
